# Adds new attendance ("asistencia") entries to the "Registro" sheet and
# tidies up the now-superseded blank "Descripcion" cells on the two most
# recent existing rows, matching the data written by the new web
# front-end's Excel-backed storage layer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The web app no longer emits an (empty) Descripcion cell for rows that
# don't have one, so drop the empty placeholders left on D10/D11.
$ws.Range("D10").ClearContents() | Out-Null
$ws.Range("D11").ClearContents() | Out-Null

# Helper: write a date-looking string (e.g. "2025-03-25") into a cell as
# literal text instead of letting Excel auto-convert it to a date serial.
# Marking the cell as Text ("@") before the assignment forces that, and
# ClearFormats() afterwards drops the number-format override again so the
# cell is left with plain default formatting, same as its neighbours.
function Set-TextValue {
    param($Cell, $Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.ClearFormats() | Out-Null
}

# New attendance rows recorded on 2025-03-25.
Set-TextValue $ws.Cells.Item(12, 1) "2025-03-25"
$ws.Cells.Item(12, 2).Value = "00:38:39"
$ws.Cells.Item(12, 3).Value = "Entrada"

Set-TextValue $ws.Cells.Item(13, 1) "2025-03-25"
$ws.Cells.Item(13, 2).Value = "00:38:42"
$ws.Cells.Item(13, 3).Value = "Entrada"

Set-TextValue $ws.Cells.Item(14, 1) "2025-03-25"
$ws.Cells.Item(14, 2).Value = "00:42:00"
$ws.Cells.Item(14, 3).Value = "Entrada"

# The very last row keeps an (empty) Descripcion cell, same shape as the
# rows before the cleanup above. Plain assignment of "" clears a cell
# outright in Excel, so use a formula that evaluates to an empty string to
# materialize an actual empty-text value in D14.
$ws.Cells.Item(14, 4).Formula = "="""""
